$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting all existing data down by one row.
$ws.Rows("1:1").Insert()

# Populate the new header row with the survey column names.
$headers = @("look", "aroma", "taste", "texture", "overall_sat", "email", "sendEmail", "gender", "age")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 8 (shifted from the old row 7) only keeps J:K as empty style cells;
# clear L8:O8 entirely so they no longer exist.
$ws.Range("L8:O8").Clear()

# Match the author's final selection on the header row.
$ws.Range("A1:I1").Select()
